# Churn Prediction Template - gender abbreviation fix
#
# The "gender" column on the first worksheet ("Feuil1") stored the male
# value using the single-letter abbreviation "M". Expand it to the full
# word "Male" so it matches the "Female, Male" convention already
# documented on the "Typology" sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A2").Value = "Male"
$ws1.Range("A3").Value = "Male"

# Leave the selection where the author left it when they saved the file.
$ws1.Range("L7").Select()
